$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.050.46"
$ws.Range("E2").Value = "  +2.10%  "
$ws.Range("D3").Value = "1.651.07"
$ws.Range("E3").Value = "  +2.09%  "
$ws.Range("E4").Value = "  +0.05%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "213.99"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.35%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.529"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("E7").Value = "  +0.02%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "23.60"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +3.28%  "
$ws.Range("E9").Value = "  +1.98%  "
$ws.Range("E10").Value = "  +0.51%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0875"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -1.46%  "
$ws.Range("D12").Value = "1.886.50"
$ws.Range("E12").Value = "  +2.21%  "
$ws.Range("D13").Value = "1.658.79"
$ws.Range("E13").Value = "  +2.56%  "
$ws.Range("E14").Value = "  +1.33%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.567"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +3.26%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "65.83"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +1.44%  "
$ws.Range("D17").Value = "28.067.48"
$ws.Range("E17").Value = "  +2.24%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "233.39"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.72%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "7.68"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("D20").Value = "0.0₃0725"
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("E21").Value = "  -0.10%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "10.72"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +5.48%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "4.41"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +2.76%  "
$ws.Range("E24").Value = "  +3.47%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "152.61"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.22%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "6.93"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +1.07%  "
$ws.Range("E27").Value = "  +0.71%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "15.80"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("E29").Value = "  +0.08%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.19"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +1.42%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.0485"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.44%  "
$ws.Range("E32").Value = "  +2.40%  "
$ws.Range("D33").Value = "1.449.51"
$ws.Range("E33").Value = "  -1.34%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "3.09"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.53%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.57"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +2.09%  "
$ws.Range("E36").Value = "  -0.38%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.894"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +4.11%  "
$ws.Range("E38").Value = "  +1.78%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.932"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -2.29%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.560"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.52%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "69.50"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +2.22%  "
$ws.Range("E42").Value = "  +3.45%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("E44").Value = "  -0.72%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.83"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +5.43%  "
$ws.Range("E46").Value = "  +3.02%  "
$ws.Range("E47").Value = "  +1.12%  "
$ws.Range("D48").Value = "1.794.81"
$ws.Range("E48").Value = "  +2.09%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "89.17"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +3.03%  "
$ws.Range("E50").Value = "  -0.57%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.101"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.55%  "
